$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2256214.8
$ws.Range("I33").Value = 3196154.2
$ws.Range("J33").Value = 360
$ws.Range("K33").Value = 3196154.2
$ws.Range("L33").Value = 360
$ws.Range("M33").Value = -3195925.2
$ws.Range("N33").Value = -818

$ws.Range("H70").Value = 2476.923
$ws.Range("I70").Value = 2400
$ws.Range("J70").Value = 2490.9092
$ws.Range("K70").Value = 7200
$ws.Range("L70").Value = 7472.7276
$ws.Range("M70").Value = -6930
$ws.Range("N70").Value = -8012.7276

$ws.Range("H73").Value = 2476.923
$ws.Range("I73").Value = 2400
$ws.Range("J73").Value = 2490.9092
$ws.Range("K73").Value = 7200
$ws.Range("L73").Value = 7472.7276
$ws.Range("M73").Value = -6264
$ws.Range("N73").Value = -9344.7276

$ws.Range("H74").Value = 12854.546
$ws.Range("I74").Value = 20000
$ws.Range("J74").Value = 4280
$ws.Range("K74").Value = 20000
$ws.Range("L74").Value = 4280
$ws.Range("M74").Value = -19064
$ws.Range("N74").Value = -6152

$ws.Range("H77").Value = 12854.546
$ws.Range("I77").Value = 20000
$ws.Range("J77").Value = 4280
$ws.Range("K77").Value = 100000
$ws.Range("L77").Value = 21400
$ws.Range("M77").Value = -95320
$ws.Range("N77").Value = -30760

$ws.Range("H112").Value = 1346.4615
$ws.Range("J112").Value = 1346.4615
$ws.Range("L112").Value = 4039.3845
$ws.Range("N112").Value = -6255.3845

$ws.Range("H129").Value = 724.5
$ws.Range("I129").Value = 485.72726
$ws.Range("J129").Value = 1600
$ws.Range("K129").Value = 1457.18178
$ws.Range("L129").Value = 4800
$ws.Range("M129").Value = 3542.81822
$ws.Range("N129").Value = -14800

$ws.Range("H138").Value = 5918.5776
$ws.Range("I138").Value = 3139.4119
$ws.Range("J138").Value = 6565.781
$ws.Range("K138").Value = 9418.235700000001
$ws.Range("L138").Value = 19697.343
$ws.Range("M138").Value = -4278.235700000001
$ws.Range("N138").Value = -29977.343

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 24084.84
$ws.Range("I32").Value = 10435.35
$ws.Range("J32").Value = 127076.45
$ws.Range("K32").Value = 10435.35
$ws.Range("L32").Value = 127076.45
$ws.Range("M32").Value = -10148.35
$ws.Range("N32").Value = -127650.45

$ws.Range("H44").Value = 25437
$ws.Range("J44").Value = 33211.8
$ws.Range("L44").Value = 33211.8
$ws.Range("N44").Value = -34187.8

$ws.Range("H74").Value = 2063.7083
$ws.Range("I74").Value = 1286.125
$ws.Range("J74").Value = 3618.875
$ws.Range("K74").Value = 1286.125
$ws.Range("L74").Value = 3618.875
$ws.Range("M74").Value = -412.125
$ws.Range("N74").Value = -5366.875

$ws.Range("H77").Value = 2063.7083
$ws.Range("I77").Value = 1286.125
$ws.Range("J77").Value = 3618.875
$ws.Range("K77").Value = 6430.625
$ws.Range("L77").Value = 18094.375
$ws.Range("M77").Value = -2062.625
$ws.Range("N77").Value = -26830.375

$ws.Range("H132").Value = 2081.362
$ws.Range("I132").Value = 1240.0377
$ws.Range("K132").Value = 3720.1131
$ws.Range("M132").Value = -1190.1131

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H25").Value = 3606
$ws.Range("I25").Value = 3606
$ws.Range("K25").Value = 3606
$ws.Range("M25").Value = -3371

$ws.Range("H107").Value = 9613.143
$ws.Range("I107").Value = 9613.143
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 9613.143
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -7693.143
$ws.Range("N107").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H110").Value = 35000
$ws.Range("J110").Value = 35000
$ws.Range("L110").Value = 35000
$ws.Range("N110").Value = -43180

$ws.Range("H111").Value = 30702
$ws.Range("J111").Value = 30702
$ws.Range("L111").Value = 30702
$ws.Range("N111").Value = -38882

$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("M115").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 750
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

$ws.Range("H25").Value = 4500
$ws.Range("I25").Value = 4500
$ws.Range("K25").Value = 4500
$ws.Range("M25").Value = -4326

$ws.Range("H31").Value = 3870.653
$ws.Range("I31").Value = 1939.1515
$ws.Range("J31").Value = 7854.375
$ws.Range("K31").Value = 1939.1515
$ws.Range("L31").Value = 7854.375
$ws.Range("M31").Value = -1644.1515
$ws.Range("N31").Value = -8444.375

$ws.Range("H34").Value = 3870.653
$ws.Range("I34").Value = 1939.1515
$ws.Range("J34").Value = 7854.375
$ws.Range("K34").Value = 1939.1515
$ws.Range("L34").Value = 7854.375
$ws.Range("M34").Value = -1737.1515
$ws.Range("N34").Value = -8258.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 71477.5
$ws.Range("I2").Value = 200056.6
$ws.Range("J2").Value = 44.666668
$ws.Range("K2").Value = 1200339.6
$ws.Range("L2").Value = 268.000008
$ws.Range("M2").Value = -1200226.6
$ws.Range("N2").Value = -494.000008

$ws.Range("H69").Value = 3088.8572
$ws.Range("I69").Value = 411
$ws.Range("J69").Value = 4160
$ws.Range("K69").Value = 1233
$ws.Range("L69").Value = 12480
$ws.Range("M69").Value = -422
$ws.Range("N69").Value = -14102

$ws.Range("H72").Value = 3088.8572
$ws.Range("I72").Value = 411
$ws.Range("J72").Value = 4160
$ws.Range("K72").Value = 3699
$ws.Range("L72").Value = 37440
$ws.Range("M72").Value = 357
$ws.Range("N72").Value = -45552

$ws.Range("H123").Value = 3141.25
$ws.Range("I123").Value = 1355
$ws.Range("J123").Value = 8500
$ws.Range("K123").Value = 4065
$ws.Range("L123").Value = 25500
$ws.Range("M123").Value = -1615
$ws.Range("N123").Value = -30400

$ws.Range("H131").Value = 889.2782999999999
$ws.Range("J131").Value = 889.2782999999999
$ws.Range("L131").Value = 2667.8349
$ws.Range("N131").Value = -12747.8349

$ws.Range("H132").Value = 886162.9399999999
$ws.Range("I132").Value = 2195149
$ws.Range("J132").Value = 13505.556
$ws.Range("K132").Value = 19756341
$ws.Range("L132").Value = 121550.004
$ws.Range("M132").Value = -19753811
$ws.Range("N132").Value = -126610.004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 5105433.5
$ws.Range("J80").Value = 5955672.5
$ws.Range("L80").Value = 5955672.5
$ws.Range("N80").Value = -5957668.5

$ws.Range("H83").Value = 5105433.5
$ws.Range("J83").Value = 5955672.5
$ws.Range("L83").Value = 29778362.5
$ws.Range("N83").Value = -29788346.5

$ws.Range("H132").Value = 2088.5
$ws.Range("I132").Value = 2155.625
$ws.Range("J132").Value = 1820
$ws.Range("K132").Value = 6466.875
$ws.Range("L132").Value = 5460
$ws.Range("M132").Value = -3936.875
$ws.Range("N132").Value = -10520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2116.6667
$ws.Range("I46").Value = 1540
$ws.Range("J46").Value = 5000
$ws.Range("K46").Value = 1540
$ws.Range("L46").Value = 5000
$ws.Range("M46").Value = -1352
$ws.Range("N46").Value = -5376

$ws.Range("H82").Value = 27780382
$ws.Range("I82").Value = 45456470
$ws.Range("J82").Value = 3671.2856
$ws.Range("K82").Value = 45456470
$ws.Range("L82").Value = 3671.2856
$ws.Range("M82").Value = -45456109
$ws.Range("N82").Value = -4393.2856

$ws.Range("H85").Value = 27780382
$ws.Range("I85").Value = 45456470
$ws.Range("J85").Value = 3671.2856
$ws.Range("K85").Value = 45456470
$ws.Range("L85").Value = 3671.2856
$ws.Range("M85").Value = -45455222
$ws.Range("N85").Value = -6167.2856

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 45005.5
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 45005.5
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 45005.5
$ws.Range("N20").Value = -45485.5
$ws.Range("M20").ClearContents()
